$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 6 (the "SOURCES OF FINANCE" table) - switch the table to a
#    different built-in table style (tableStyleId).
# ---------------------------------------------------------------------------
$tableSlide = $p.Slides.Item(6)
$tableShape = $tableSlide.Shapes.Item(2)
$tbl = $tableShape.Table
$tbl.ApplyStyle("{66E45740-D788-412D-A827-DDFFA8C2B17D}")

# ---------------------------------------------------------------------------
# 2) Theme colours - the deck's main theme ("Integral") and the notes
#    theme ("Office Theme") were swapped. The slide-facing colour scheme
#    (theme1.xml, reached through any slide's ThemeColorScheme) becomes the
#    plain default "Office Theme" palette.
# ---------------------------------------------------------------------------
$officeThemeColors = @(
    0x000000,  # 1  dk1
    0xFFFFFF,  # 2  lt1
    0x44546A,  # 3  dk2
    0xE7E6E6,  # 4  lt2
    0x5B9BD5,  # 5  accent1
    0xED7D31,  # 6  accent2
    0xA5A5A5,  # 7  accent3
    0xFFC000,  # 8  accent4
    0x4472C4,  # 9  accent5
    0x70AD47,  # 10 accent6
    0x0563C1,  # 11 hlink
    0x954F72   # 12 folHlink
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $rgbHex = $officeThemeColors[$i - 1]
    $r = $rgbHex -band 0xFF
    $g = ($rgbHex -shr 8) -band 0xFF
    $b = ($rgbHex -shr 16) -band 0xFF
    $comRgb = $r + ($g * 256) + ($b * 65536)
    $tcs.Colors($i).RGB = $comRgb
}
